$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.116.90"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "3.058.18"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "386.95"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "101.99"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").Value = "0.536"
$ws.Range("E7").Value = "  -1.91%  "
$ws.Range("D9").Value = "0.579"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "36.73"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "0.0847"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").Value = "3.535.45"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "18.29"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").Value = "7.68"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "3.057.08"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "0.983"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "10.66"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "51.100.94"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").Value = "3.20"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").Value = "0.0₃0953"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("D22").Value = "12.23"
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("D23").Value = "69.53"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "264.02"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "3.12"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "7.89"
$ws.Range("E26").Value = "  -5.14%  "
$ws.Range("D27").Value = "26.99"
$ws.Range("E27").Value = "  +2.98%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "7.21"
$ws.Range("E28").Value = "  -4.16%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "0.162"
$ws.Range("E30").Value = "  -5.86%  "
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("D32").Value = "10.35"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "35.55"
$ws.Range("E33").Value = "  +4.64%  "
$ws.Range("D34").Value = "0.0471"
$ws.Range("E34").Value = "  +4.32%  "
$ws.Range("D35").Value = "2.07"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "50.03"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "3.36"
$ws.Range("E38").Value = "  +2.06%  "
$ws.Range("D39").Value = "0.295"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "130.20"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("D41").Value = "1.82"
$ws.Range("E41").Value = "  -1.59%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "16.48"
$ws.Range("E42").Value = "  -2.87%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "0.115"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").Value = "3.75"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("D46").Value = "21.59"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "2.50"
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "2.051.90"
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("D50").Value = "0.0324"
$ws.Range("E50").Value = "  +1.70%  "
$ws.Range("D51").Value = "0.903"
$ws.Range("E51").Value = "  +14.50%  "
